$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 46
$ws.Range("A46").Value = 43204
$ws.Range("B46").Value = "Préparation de la présentation intermédiaire et présentation"
$ws.Range("C46").Value = 1.5

# Row 47
$ws.Range("A47").Value = 43205
$ws.Range("B47").Value = "Clean du repo git, suppression des fichiers inutilisés"
$ws.Range("C47").Value = 0.25

# Row 48
$ws.Range("A48").Value = 43208
$ws.Range("B48").Value = "Implémentation des classes CategoryModel et CategoryLogic. Adaptation des classes précédemment faites."
$ws.Range("C48").Value = 1.5
$ws.Rows.Item(48).RowHeight = 30

# Update selection to match the diff
$ws.Range("A49").Select()
